$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D, shifting existing quarters from D:K to F:M
$ws.Columns("D:E").Insert()

# Copy formats (number/date styles) from column F into new D:E for rows 7-102
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Rows 36 and 78 are fully blank separator rows, and rows 37/79 are section title rows
# with no data columns - remove the stray formatted cells the paste created there.
$ws.Range("D36:E36").Clear()
$ws.Range("D37:E37").Clear()
$ws.Range("D78:E78").Clear()
$ws.Range("D79:E79").Clear()

# Populate the two new quarter columns (D = period 43496, E = period 43404)
$ws.Range("D7").Value = 43496
$ws.Range("E7").Value = 43404
$ws.Range("D8").Value = 384100
$ws.Range("E8").Value = 424900
$ws.Range("D9").Value = 307200
$ws.Range("E9").Value = 338100
$ws.Range("D10").Value = 76900
$ws.Range("E10").Value = 86800
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = -400
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 356900
$ws.Range("E17").Value = 389400
$ws.Range("D18").Value = 27200
$ws.Range("E18").Value = 35500
$ws.Range("D20").Value = -3000
$ws.Range("E20").Value = -10100
$ws.Range("D21").Value = 47700
$ws.Range("E21").Value = 49100
$ws.Range("D22").Value = "NA"
$ws.Range("E22").Value = "NA"
$ws.Range("D23").Value = 24100
$ws.Range("E23").Value = 25400
$ws.Range("D24").Value = 5700
$ws.Range("E24").Value = 6900
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 18400
$ws.Range("E26").Value = 18500
$ws.Range("D27").Value = 18400
$ws.Range("E27").Value = 18500
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 3000
$ws.Range("E32").Value = 10100
$ws.Range("D33").Value = 18400
$ws.Range("E33").Value = 18500
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 18400
$ws.Range("E35").Value = 18500
$ws.Range("D38").Value = 43496
$ws.Range("E38").Value = 43404
$ws.Range("D41").Value = 42000
$ws.Range("E41").Value = 57900
$ws.Range("D42").Value = 2500
$ws.Range("E42").Value = 4500
$ws.Range("D43").Value = 118000
$ws.Range("E43").Value = 136500
$ws.Range("D44").Value = 116100
$ws.Range("E44").Value = 116000
$ws.Range("D45").Value = 13900
$ws.Range("E45").Value = 11700
$ws.Range("D46").Value = 292500
$ws.Range("E46").Value = 326500
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 212000
$ws.Range("E48").Value = 213400
$ws.Range("D49").Value = 998500
$ws.Range("E49").Value = 1010800
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 29800
$ws.Range("E52").Value = 26400
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 1532800
$ws.Range("E54").Value = 1577100
$ws.Range("D57").Value = 57000
$ws.Range("E57").Value = 66200
$ws.Range("D58").Value = 2300
$ws.Range("E58").Value = 4400
$ws.Range("D59").Value = 83200
$ws.Range("E59").Value = 95600
$ws.Range("D60").Value = 142500
$ws.Range("E60").Value = 166300
$ws.Range("D61").Value = 709800
$ws.Range("E61").Value = 717900
$ws.Range("D62").Value = 72500
$ws.Range("E62").Value = 76600
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 924900
$ws.Range("E66").Value = 960800
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 302200
$ws.Range("E72").Value = 303200
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 608000
$ws.Range("E76").Value = 616300
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43496
$ws.Range("E80").Value = 43404
$ws.Range("D81").Value = 18400
$ws.Range("E81").Value = 18500
$ws.Range("D83").Value = 23600
$ws.Range("E83").Value = 23700
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 30300
$ws.Range("E89").Value = 54700
$ws.Range("D91").Value = -11500
$ws.Range("E91").Value = -5400
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -11600
$ws.Range("E94").Value = -3300
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -34600
$ws.Range("E100").Value = -43700
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = -15900
$ws.Range("E102").Value = 7700

# Fix one historical data correction uncovered by re-shifting: column I row 14 becomes 0 (was NA)
$ws.Range("I14").Value = 0
